$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'43.017.89"
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = "'  -0.35%  "
$ws.Range('E2').Style = 'Normal'
$ws.Range('D3').Value = "'2.303.67"
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = "'  -0.53%  "
$ws.Range('E3').Style = 'Normal'
$ws.Range('E4').Value = "'  +0.02%  "
$ws.Range('E4').Style = 'Normal'
$ws.Range('D5').Value = "'301.84"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = "'  -0.43%  "
$ws.Range('E5').Style = 'Normal'
$ws.Range('D6').Value = "'98.63"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = "'  -3.12%  "
$ws.Range('E6').Style = 'Normal'
$ws.Range('D7').Value = "'0.525"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = "'  +3.86%  "
$ws.Range('E7').Style = 'Normal'
$ws.Range('E8').Value = "'  +0.03%  "
$ws.Range('E8').Style = 'Normal'
$ws.Range('D9').Value = "'0.523"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = "'  +0.61%  "
$ws.Range('E9').Style = 'Normal'
$ws.Range('D10').Value = "'35.70"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = "'  -1.11%  "
$ws.Range('E10').Style = 'Normal'
$ws.Range('D11').Value = "'0.0790"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = "'  -0.82%  "
$ws.Range('E11').Style = 'Normal'
$ws.Range('E12').Value = "'  -1.22%  "
$ws.Range('E12').Style = 'Normal'
$ws.Range('D13').Value = "'17.93"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = "'  +0.08%  "
$ws.Range('E13').Style = 'Normal'
$ws.Range('D14').Value = "'6.89"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = "'  -0.19%  "
$ws.Range('E14').Style = 'Normal'
$ws.Range('D15').Value = "'2.662.48"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = "'  -0.23%  "
$ws.Range('E15').Style = 'Normal'
$ws.Range('D16').Value = "'2.267.02"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = "'  -1.80%  "
$ws.Range('E16').Style = 'Normal'
$ws.Range('D17').Value = "'0.789"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = "'  -3.13%  "
$ws.Range('E17').Style = 'Normal'
$ws.Range('D18').Value = "'42.905.92"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = "'  -0.38%  "
$ws.Range('E18').Style = 'Normal'
$ws.Range('D19').Value = "'13.39"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = "'  +5.95%  "
$ws.Range('E19').Style = 'Normal'
$ws.Range('D20').Value = "'0.0₃0909"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = "'  +0.45%  "
$ws.Range('E20').Style = 'Normal'
$ws.Range('D21').Value = "'6.18"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = "'  -0.01%  "
$ws.Range('E21').Style = 'Normal'
$ws.Range('D22').Value = "'68.25"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = "'  +0.50%  "
$ws.Range('E22').Style = 'Normal'
$ws.Range('D23').Value = "'239.56"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = "'  +0.85%  "
$ws.Range('E23').Style = 'Normal'
$ws.Range('E24').Value = "'  -3.03%  "
$ws.Range('E24').Style = 'Normal'
$ws.Range('B25').Value = "'Dai"
$ws.Range('B25').Style = 'Normal'
$ws.Range('C25').Value = "'https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range('C25').Style = 'Normal'
$ws.Range('D25').Value = "'0.999"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = "'  -0.11%  "
$ws.Range('E25').Style = 'Normal'
$ws.Range('B26').Value = "'PancakeSwap"
$ws.Range('B26').Style = 'Normal'
$ws.Range('C26').Value = "'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range('C26').Style = 'Normal'
$ws.Range('D26').Value = "'2.44"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = "'  -1.51%  "
$ws.Range('E26').Style = 'Normal'
$ws.Range('D27').Value = "'24.78"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = "'  -0.32%  "
$ws.Range('E27').Style = 'Normal'
$ws.Range('D28').Value = "'168.12"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = "'  -0.02%  "
$ws.Range('E28').Style = 'Normal'
$ws.Range('B29').Value = "'Toncoin"
$ws.Range('B29').Style = 'Normal'
$ws.Range('C29').Value = "'https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range('C29').Style = 'Normal'
$ws.Range('D29').Value = "'2.04"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = "'  -12.47%  "
$ws.Range('E29').Style = 'Normal'
$ws.Range('B30').Value = "'Cosmos"
$ws.Range('B30').Style = 'Normal'
$ws.Range('C30').Value = "'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range('C30').Style = 'Normal'
$ws.Range('D30').Value = "'9.14"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = "'  -1.16%  "
$ws.Range('E30').Style = 'Normal'
$ws.Range('D31').Value = "'33.35"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = "'  -3.46%  "
$ws.Range('E31').Style = 'Normal'
$ws.Range('D32').Value = "'5.21"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = "'  +3.24%  "
$ws.Range('E32').Style = 'Normal'
$ws.Range('D34').Value = "'4.84"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = "'  +1.63%  "
$ws.Range('E34').Style = 'Normal'
$ws.Range('D35').Value = "'18.23"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = "'  +5.18%  "
$ws.Range('E35').Style = 'Normal'
$ws.Range('D36').Value = "'2.41"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = "'  -0.57%  "
$ws.Range('E36').Style = 'Normal'
$ws.Range('E37').Value = "'  -0.51%  "
$ws.Range('E37').Style = 'Normal'
$ws.Range('B38').Value = "'ARBITRUM"
$ws.Range('B38').Style = 'Normal'
$ws.Range('C38').Value = "'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range('C38').Style = 'Normal'
$ws.Range('D38').Value = "'1.79"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = "'  -0.20%  "
$ws.Range('E38').Style = 'Normal'
$ws.Range('B39').Value = "'Kaspa"
$ws.Range('B39').Style = 'Normal'
$ws.Range('C39').Value = "'https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range('C39').Style = 'Normal'
$ws.Range('D39').Value = "'0.101"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = "'  -1.51%  "
$ws.Range('E39').Style = 'Normal'
$ws.Range('E40').Value = "'  +1.76%  "
$ws.Range('E40').Style = 'Normal'
$ws.Range('E41').Value = "'  -2.95%  "
$ws.Range('E41').Style = 'Normal'
$ws.Range('D42').Value = "'1.998.08"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = "'  +0.20%  "
$ws.Range('E42').Style = 'Normal'
$ws.Range('E43').Value = "'  -0.49%  "
$ws.Range('E43').Style = 'Normal'
$ws.Range('B44').Value = "'FraxShare"
$ws.Range('B44').Style = 'Normal'
$ws.Range('C44').Value = "'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range('C44').Style = 'Normal'
$ws.Range('D44').Value = "'10.08"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = "'  -1.48%  "
$ws.Range('E44').Style = 'Normal'
$ws.Range('B45').Value = "'ApeXProtocol"
$ws.Range('B45').Style = 'Normal'
$ws.Range('C45').Value = "'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range('C45').Style = 'Normal'
$ws.Range('D45').Value = "'2.12"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = "'  -7.79%  "
$ws.Range('E45').Style = 'Normal'
$ws.Range('D46').Value = "'17.47"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = "'  -1.06%  "
$ws.Range('E46').Style = 'Normal'
$ws.Range('E47').Value = "'  -2.69%  "
$ws.Range('E47').Style = 'Normal'
$ws.Range('E48').Value = "'  -2.32%  "
$ws.Range('E48').Style = 'Normal'
$ws.Range('D49').Value = "'2.531.34"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = "'  -0.67%  "
$ws.Range('E49').Style = 'Normal'
$ws.Range('E50').Value = "'  -0.12%  "
$ws.Range('E50').Style = 'Normal'
$ws.Range('D51').Value = "'73.36"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = "'  +4.16%  "
$ws.Range('E51').Style = 'Normal'
